# Auto-generated edit script: apply cryptos price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D that hold numeric-looking text must be pinned to Text format first,
# otherwise Excel auto-converts strings like "4.40" into the number 4.4 and the
# original text representation (trailing zeros, etc.) is lost.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply the updated cell values (prices, 1h volume deltas, and the two swapped
# coin rows 44/45 with their name/link/price/volume).
$ws.Range('D2').Value = '26.843.64'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '1.640.07'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('E4').Value = '  -0.63%  '
$ws.Range('D5').Value = '217.03'
$ws.Range('E5').Value = '  -0.78%  '
$ws.Range('D6').Value = '0.514'
$ws.Range('E6').Value = '  +2.87%  '
$ws.Range('E7').Value = '  -0.58%  '
$ws.Range('D8').Value = '0.255'
$ws.Range('E8').Value = '  +1.84%  '
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('D10').Value = '19.93'
$ws.Range('E10').Value = '  +3.77%  '
$ws.Range('D11').Value = '0.0845'
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').Value = '1.869.10'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').Value = '1.640.55'
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('D14').Value = '4.12'
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('D15').Value = '0.530'
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('D16').Value = '66.83'
$ws.Range('E16').Value = '  +3.38%  '
$ws.Range('D17').Value = '26.841.91'
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').Value = '0.0₃0728'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').Value = '219.62'
$ws.Range('E19').Value = '  +2.34%  '
$ws.Range('E20').Value = '  -0.57%  '
$ws.Range('D21').Value = '6.74'
$ws.Range('E21').Value = '  +7.01%  '
$ws.Range('D22').Value = '4.40'
$ws.Range('E22').Value = '  +1.02%  '
$ws.Range('E23').Value = '  +3.56%  '
$ws.Range('E24').Value = '  +0.72%  '
$ws.Range('D25').Value = '147.56'
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('E27').Value = '  +4.87%  '
$ws.Range('E28').Value = '  +1.10%  '
$ws.Range('D29').Value = '15.80'
$ws.Range('E29').Value = '  +0.78%  '
$ws.Range('D30').Value = '0.0504'
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  -1.47%  '
$ws.Range('D32').Value = '3.34'
$ws.Range('E32').Value = '  -1.24%  '
$ws.Range('E33').Value = '  +0.85%  '
$ws.Range('E34').Value = '  +1.23%  '
$ws.Range('D35').Value = '1.261.49'
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('D36').Value = '2.43'
$ws.Range('E36').Value = '  -0.49%  '
$ws.Range('E37').Value = '  +1.91%  '
$ws.Range('E38').Value = '  +1.28%  '
$ws.Range('D39').Value = '0.833'
$ws.Range('E39').Value = '  +2.46%  '
$ws.Range('E40').Value = '  -0.54%  '
$ws.Range('D41').Value = '0.808'
$ws.Range('E41').Value = '  +0.32%  '
$ws.Range('E42').Value = '  +2.57%  '
$ws.Range('D43').Value = '1.782.83'
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = '2.10'
$ws.Range('E44').Value = '  -1.45%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '61.74'
$ws.Range('E45').Value = '  +2.86%  '
$ws.Range('D46').Value = '91.92'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').Value = '1.58'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('E48').Value = '  +17.60%  '
$ws.Range('D49').Value = '0.0514'
$ws.Range('E49').Value = '  -0.47%  '
$ws.Range('D50').Value = '7.64'
$ws.Range('E50').Value = '  +1.94%  '
$ws.Range('D51').Value = '0.0965'
$ws.Range('E51').Value = '  +0.38%  '

